$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.75"
$ws.Range("D3").Value = "'25.06"
$ws.Range("D4").Value = "'5.151"
$ws.Range("D5").Value = "'0.05637"
$ws.Range("D6").Value = "'6.518"
$ws.Range("D7").Value = "'2.982"
$ws.Range("D8").Value = "'0.8134"
$ws.Range("D9").Value = "'0.8364"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.009520"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1330"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.06955"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02839"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09402"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001517"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("D16").Value = "'0.006270"
$ws.Range("D17").Value = "'3.504"
$ws.Range("D18").Value = "'2.107"
$ws.Range("D19").Value = "'0.3166"
$ws.Range("D20").Value = "'0.03313"
$ws.Range("D21").Value = "'0.1319"
$ws.Range("D22").Value = "'3.744"
$ws.Range("D23").Value = "'0.04703"
$ws.Range("D24").Value = "'0.1370"
$ws.Range("D25").Value = "'0.001239"
$ws.Range("D26").Value = "'0.004529"
$ws.Range("D27").Value = "'0.00009701"
$ws.Range("E27").Value = "26NitroExNTX"
$ws.Range("D28").Value = "'0.0001940"
$ws.Range("D40").Value = "'0.03624"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006270"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1050"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002717"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.007589"
$ws.Range("D45").Value = "'0.00005290"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.2200"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "'0.002286"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
